# Mission 3, 4, and 5 notes added to DUNEXMainExp notes workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5 (Mission 3), Row 6 (Mission 4), Row 7 (Mission 5) -----------------
# Columns A-E, then F (new shared string "small grid (3x3)"), then G-M,
# all written row by row first (these columns only reuse already-existing
# shared strings so ordering among them does not affect the shared string table).

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "0.5 - 1"
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = "E - SE"
$ws.Range("E5").Value = "surf board"
$ws.Range("F5").Value = "small grid (3x3)"
$ws.Range("G5").Value = "Alex de Klerk"
$ws.Range("H5").Value = "Sean McGill"
$ws.Range("I5").Value = "Christine Baker"
$ws.Range("J5").Value = "EJ Rainville"
$ws.Range("K5").Value = "EJ Rainville"
$ws.Range("L5").Value = 9
$ws.Range("M5").Value = 2

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "0.5 - 1"
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = "E - SE"
$ws.Range("E6").Value = "surf board"
$ws.Range("F6").Value = "small grid (3x3)"
$ws.Range("G6").Value = "Alex de Klerk"
$ws.Range("H6").Value = "Sean McGill"
$ws.Range("I6").Value = "Christine Baker"
$ws.Range("J6").Value = "EJ Rainville"
$ws.Range("K6").Value = "EJ Rainville"
$ws.Range("L6").Value = 9
$ws.Range("M6").Value = 2

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "0.5 - 1"
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = "E - SE"
$ws.Range("E7").Value = "surf board"
$ws.Range("F7").Value = "small grid (3x3)"
$ws.Range("G7").Value = "Alex de Klerk"
$ws.Range("H7").Value = "Sean McGill"
$ws.Range("I7").Value = "Christine Baker"
$ws.Range("J7").Value = "EJ Rainville"
$ws.Range("K7").Value = "EJ Rainville"
$ws.Range("L7").Value = 9
$ws.Range("M7").Value = 2

# --- Start / End times (new shared strings) ---------------------------------
$ws.Range("Q5").Value = "2021-10-05T17:10:00"
$ws.Range("R5").Value = "2021-10-05T17:38:00"
$ws.Range("Q6").Value = "2021-10-05T18:06:00"
$ws.Range("R6").Value = "2021-10-05T18:25:00"
$ws.Range("Q7").Value = "2021-10-05T18:35:00"
$ws.Range("R7").Value = "2021-10-05T18:48:00"

# --- Deployment Notes (new shared strings) -----------------------------------
$ws.Range("S5").Value = "All microSWIFTs were retrieved" + [char]10 + "- We needed to play goalie at the pier since the micros were headed north quickly - this was effective for getting them all back" + [char]10 + "- this mission was right at low tide"
$ws.Range("S6").Value = "This mission was just after low tide and there was a light onshore breeze"
$ws.Range("S7").Value = "This mission was very short "

# --- microSWIFTs Deployed / Retrieved / Shepherds Retrieved -------------------
$ws.Range("N5").Value = "3,4,5,40,41,42,39,38,57"
$ws.Range("O5").Value = "3,4,5,40,41,42,39,38,57"
$ws.Range("P5").Value = "3,4"

$ws.Range("N6").Value = "3,4,5,40,41,42,39,38,57"
$ws.Range("O6").Value = "3,4,5,40,41,42,39,38,57"
$ws.Range("P6").Value = "3,4"

$ws.Range("N7").Value = "3,4,5,40,41,42,39,38,57"
$ws.Range("O7").Value = "3,4,5,40,41,42,39,38,57"
$ws.Range("P7").Value = "3,4"

# --- Column S autofit (was widened from the long note text) ------------------
$ws.Columns.Item(19).ColumnWidth = 108

# --- Row heights (autofit for the new wrapped-text rows) ---------------------
$ws.Rows.Item(5).RowHeight = 51
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 17

# --- Final selection / scroll position ---------------------------------------
$ws.Range("L1").Select()
$ws.Range("O7").Select()
